$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header C1: "audioFalse" -> "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# C2: "trainingaudio/26_kapako1.wav" -> "train1P2"
$ws.Range("C2").Value = "train1P2"

# C3: "trainingaudio/08_tipako2.wav" -> "train1P2"
$ws.Range("C3").Value = "train1P2"
